$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7774.75
$ws.Range("I19").Value = 12199.75
$ws.Range("J19").Value = 3349.75
$ws.Range("K19").Value = 12199.75
$ws.Range("L19").Value = 3349.75
$ws.Range("M19").Value = -12024.75
$ws.Range("N19").Value = -3699.75

$ws.Range("H75").Value = 16933.334
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872

$ws.Range("H78").Value = 16933.334
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1084.5883
$ws.Range("I74").Value = 1036.5
$ws.Range("K74").Value = 1036.5
$ws.Range("M74").Value = -162.5

$ws.Range("H77").Value = 1084.5883
$ws.Range("I77").Value = 1036.5
$ws.Range("K77").Value = 5182.5
$ws.Range("M77").Value = -814.5

$ws.Range("H132").Value = 5963.731
$ws.Range("I132").Value = 7496.8823
$ws.Range("J132").Value = 3067.7778
$ws.Range("K132").Value = 22490.6469
$ws.Range("L132").Value = 9203.3334
$ws.Range("M132").Value = -19960.6469
$ws.Range("N132").Value = -14263.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 43093.92
$ws.Range("I134").Value = 74503.43
$ws.Range("J134").Value = 3118.182
$ws.Range("K134").Value = 223510.29
$ws.Range("L134").Value = 9354.545999999998
$ws.Range("M134").Value = -220975.29
$ws.Range("N134").Value = -14424.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3502
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12530
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 817.73334
$ws.Range("I5").Value = 472.42856
$ws.Range("J5").Value = 1119.875
$ws.Range("K5").Value = 1417.28568
$ws.Range("L5").Value = 3359.625
$ws.Range("M5").Value = -1305.28568
$ws.Range("N5").Value = -3583.625

$ws.Range("H122").Value = 1460.4762
$ws.Range("I122").Value = 1610.3334
$ws.Range("J122").Value = 1348.0834
$ws.Range("K122").Value = 14493.0006
$ws.Range("L122").Value = 12132.7506
$ws.Range("M122").Value = -12043.0006
$ws.Range("N122").Value = -17032.7506

$ws.Range("H132").Value = 2141.3333
$ws.Range("I132").Value = 1148
$ws.Range("J132").Value = 2936
$ws.Range("K132").Value = 10332
$ws.Range("L132").Value = 26424
$ws.Range("M132").Value = -7802
$ws.Range("N132").Value = -31484

$ws.Range("H135").Value = 817.73334
$ws.Range("I135").Value = 472.42856
$ws.Range("J135").Value = 1119.875
$ws.Range("K135").Value = 4251.85704
$ws.Range("L135").Value = 10078.875
$ws.Range("M135").Value = -1716.85704
$ws.Range("N135").Value = -15148.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 73.07692
$ws.Range("I2").Value = 31.8
$ws.Range("K2").Value = 31.8
$ws.Range("M2").Value = 81.2

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25001366
$ws.Range("J7").Value = 71430820
$ws.Range("L7").Value = 71430820
$ws.Range("N7").Value = -71431044

$ws.Range("H40").Value = 3250
$ws.Range("I40").Value = 3500
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3364
$ws.Range("N40").Value = -3272

$ws.Range("H124").Value = 46750
$ws.Range("J124").Value = 46750
$ws.Range("L124").Value = 46750
$ws.Range("N124").Value = -56570

$ws.Range("H126").Value = 25001366
$ws.Range("J126").Value = 71430820
$ws.Range("L126").Value = 214292460
$ws.Range("N126").Value = -214297400

$ws.Range("H136").Value = 8168.625
$ws.Range("I136").Value = 11856.444
$ws.Range("J136").Value = 3427.1428
$ws.Range("K136").Value = 35569.33199999999
$ws.Range("L136").Value = 10281.4284
$ws.Range("M136").Value = -33019.33199999999
$ws.Range("N136").Value = -15381.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 20692.334
$ws.Range("J80").Value = 20692.334
$ws.Range("L80").Value = 20692.334
$ws.Range("N80").Value = -22688.334

$ws.Range("H83").Value = 20692.334
$ws.Range("J83").Value = 20692.334
$ws.Range("L83").Value = 62077.00199999999
$ws.Range("N83").Value = -72061.002

$ws.Range("H122").Value = 170810.5
$ws.Range("I122").Value = 338334.66
$ws.Range("J122").Value = 3286.3333
$ws.Range("K122").Value = 1015003.98
$ws.Range("L122").Value = 9858.999899999999
$ws.Range("M122").Value = -1012553.98
$ws.Range("N122").Value = -14758.9999

$ws.Range("H124").Value = 44714.5
$ws.Range("J124").Value = 44714.5
$ws.Range("L124").Value = 44714.5
$ws.Range("N124").Value = -54534.5

$ws.Range("H125").Value = 40238.332
$ws.Range("J125").Value = 40238.332
$ws.Range("L125").Value = 40238.332
$ws.Range("N125").Value = -50078.332

$ws.Range("H126").Value = 19534
$ws.Range("I126").Value = 22840.8
$ws.Range("K126").Value = 68522.4
$ws.Range("M126").Value = -66052.4

$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920

$ws.Range("H136").Value = 6852.0835
$ws.Range("I136").Value = 8091.6665
$ws.Range("J136").Value = 3133.3333
$ws.Range("K136").Value = 24274.9995
$ws.Range("L136").Value = 9399.999899999999
$ws.Range("M136").Value = -21724.9995
$ws.Range("N136").Value = -14499.9999
